# Update cryptos list values (prices and 1h volume change) per the
# scraped-data refresh commit.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '67.678.88'
$ws.Range("E2").Value = '  -0.73%  '
$ws.Range("D3").Value = '3.727.04'
$ws.Range("E3").Value = '  -2.84%  '
$ws.Range("D4").Value = '''1.00'
$ws.Range("E4").Value = '  +0.28%  '
$ws.Range("D5").Value = '''590.59'
$ws.Range("E5").Value = '  -1.59%  '
$ws.Range("D6").Value = '''164.80'
$ws.Range("E6").Value = '  -3.93%  '
$ws.Range("D7").Value = '3.729.50'
$ws.Range("E7").Value = '  -2.77%  '
$ws.Range("E8").Value = '  +0.05%  '
$ws.Range("D9").Value = '''0.516'
$ws.Range("E9").Value = '  -1.73%  '
$ws.Range("D10").Value = '''0.158'
$ws.Range("E10").Value = '  -5.40%  '
$ws.Range("D11").Value = '''6.39'
$ws.Range("E11").Value = '  -1.76%  '
$ws.Range("D12").Value = '''0.445'
$ws.Range("E12").Value = '  -3.36%  '
$ws.Range("D13").Value = '''0.0000260'
$ws.Range("E13").Value = '  -8.48%  '
$ws.Range("D14").Value = '''35.59'
$ws.Range("E14").Value = '  -3.60%  '
$ws.Range("D15").Value = '4.375.01'
$ws.Range("E15").Value = '  -2.28%  '
$ws.Range("D16").Value = '3.746.51'
$ws.Range("E16").Value = '  -2.50%  '
$ws.Range("D17").Value = '67.807.51'
$ws.Range("E17").Value = '  -0.66%  '
$ws.Range("D18").Value = '''18.34'
$ws.Range("E18").Value = '  -0.09%  '
$ws.Range("D19").Value = '''6.97'
$ws.Range("E19").Value = '  -6.22%  '
$ws.Range("E20").Value = '  -0.28%  '
$ws.Range("D21").Value = '''10.40'
$ws.Range("E21").Value = '  -4.29%  '
$ws.Range("D22").Value = '''461.65'
$ws.Range("E22").Value = '  -1.37%  '
$ws.Range("D23").Value = '''0.695'
$ws.Range("E23").Value = '  -4.81%  '
$ws.Range("D24").Value = '''82.58'
$ws.Range("E24").Value = '  -1.04%  '
$ws.Range("D25").Value = '''0.0000133'
$ws.Range("E25").Value = '  -15.82%  '
$ws.Range("D26").Value = '''2.15'
$ws.Range("E26").Value = '  -5.36%  '
$ws.Range("D27").Value = '''11.84'
$ws.Range("E27").Value = '  -2.34%  '
$ws.Range("D28").Value = '''10.05'
$ws.Range("E28").Value = '  -3.49%  '
$ws.Range("E29").Value = '  +0.00%  '
$ws.Range("D30").Value = '3.892.83'
$ws.Range("E30").Value = '  -2.41%  '
$ws.Range("D31").Value = '''2.86'
$ws.Range("E31").Value = '  -2.41%  '
$ws.Range("D32").Value = '''7.28'
$ws.Range("E32").Value = '  -5.95%  '
$ws.Range("D33").Value = '''29.53'
$ws.Range("E33").Value = '  -4.73%  '
$ws.Range("D34").Value = '''2.16'
$ws.Range("E34").Value = '  -5.67%  '
$ws.Range("D35").Value = '''8.93'
$ws.Range("E35").Value = '  -4.84%  '
$ws.Range("D36").Value = '3.694.96'
$ws.Range("E36").Value = '  -2.73%  '
$ws.Range("D37").Value = '''0.101'
$ws.Range("E37").Value = '  -4.14%  '
$ws.Range("D38").Value = '''3.42'
$ws.Range("E38").Value = '  -11.97%  '
$ws.Range("B39").Value = 'Kaspa'
$ws.Range("C39").Value = 'https://coinranking.com/coin/V8GxkwWow+kaspa-kas'
$ws.Range("D39").Value = '''0.136'
$ws.Range("E39").Value = '  -2.29%  '
$ws.Range("B40").Value = 'Mantle'
$ws.Range("C40").Value = 'https://coinranking.com/coin/BoI4ux0nd+mantle-mnt'
$ws.Range("D40").Value = '''0.989'
$ws.Range("E40").Value = '  -2.86%  '
$ws.Range("D41").Value = '''5.70'
$ws.Range("E41").Value = '  -4.52%  '
$ws.Range("D42").Value = '''1.00'
$ws.Range("E42").Value = '  +0.07%  '
$ws.Range("D44").Value = '''0.303'
$ws.Range("E44").Value = '  -4.72%  '
$ws.Range("D45").Value = '''8.47'
$ws.Range("E45").Value = '  -3.34%  '
$ws.Range("D46").Value = '''1.90'
$ws.Range("E46").Value = '  -4.22%  '
$ws.Range("D47").Value = '''45.27'
$ws.Range("E47").Value = '  -2.80%  '
$ws.Range("D48").Value = '''390.11'
$ws.Range("E48").Value = '  -6.53%  '
$ws.Range("D49").Value = '''144.39'
$ws.Range("E49").Value = '  +2.02%  '
$ws.Range("D50").Value = '''0.0343'
$ws.Range("E50").Value = '  -4.54%  '
$ws.Range("D51").Value = '''37.95'
$ws.Range("E51").Value = '  -0.56%  '
